$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.956.56'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -2.58%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.861.39'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -2.01%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '305.92'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.19%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5065'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.98%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3731'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.99%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07132'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8872'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.52'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.99%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07548'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.89%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.857.34'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -2.84%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.289'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.65%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '88.93'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -2.81%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.002'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008366'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -4.54%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.05'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -2.80%  '
$ws.Range('E19').Value = '  +0.18%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '26.999.89'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -2.54%  '
$ws.Range('E21').Value = '  -1.75%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.090.63'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.49'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.81%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.461'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.91%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.846'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '147.33'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -3.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.95'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.64%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.088'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -4.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '112.58'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.74%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.649'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -3.56%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.643'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -2.91%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09046'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.37%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05112'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -3.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.050'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -4.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.148'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -6.57%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7271'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -7.30%  '
$ws.Range('E37').Value = '  -2.20%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.039'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.39%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.444'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -6.73%  '
$ws.Range('E40').Value = '  -1.72%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5311'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -3.14%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.574'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.15%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '115.29'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.269'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.91%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1469'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -2.49%  '
$ws.Range('E46').Value = '  +0.25%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4600'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -3.47%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.998'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -4.14%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.555'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -3.48%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '36.47'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.71%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '63.88'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -4.03%  '
